$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Paragraph 5 currently reads "Also selector methods".
$para5 = $tr.Paragraphs(5, 1)

# Insert the new sentence (plus a paragraph break) before the existing
# "Also selector methods" text. This pushes "Also selector methods" down
# into a new paragraph 6, while the new sentence becomes paragraph 5 -
# inheriting the Wingdings symbol-font run formatting already present on
# that paragraph.
$null = $para5.InsertBefore("Returns a COPY of the element, not a reference`r")

# Re-fetch paragraph 5 (the new sentence) and recolor just "COPY" in red,
# matching the existing red accent runs used elsewhere on this slide.
$tr = $shape.TextFrame.TextRange
$newPara5 = $tr.Paragraphs(5, 1)
$copyRange = $newPara5.Characters(11, 4)
$copyRange.Font.Color.RGB = 255

# Re-fetch paragraph 6 ("Also selector methods") and split it into two
# runs - "Also " and "selector methods" - without altering formatting.
# Re-assigning the same text to a sub-range forces PowerPoint to break
# the run in two while leaving the inherited run formatting untouched.
$tr = $shape.TextFrame.TextRange
$para6 = $tr.Paragraphs(6, 1)
$alsoRange = $para6.Characters(1, 5)
$alsoRange.Text = "Also "
